# "Fixed issue where site would not record"
#
# 1. Rounds sheet: three rounds (rows 8, 9, 10 for Match ID row "9") were
#    missing their "Site" value (column D) - backfill them.
#    Also round 83's Attack/Defense (column E) was mis-recorded as Attack,
#    should be Defense.
# 2. Players sheet: rounds 7-9 of match "9" (rows 32-46) recorded the map
#    name in the "Spawn" column (E) instead of the actual bombsite - fix
#    them to the correct site strings.
# 3. Operator names (Players!D) were stored title-cased; normalize them to
#    lowercase (and fix "Jager" -> "jäger" with the umlaut) everywhere they
#    appear.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rounds sheet fixes
# ---------------------------------------------------------------------
$rounds = $wb.Worksheets.Item("Rounds")

$rounds.Range("D8").Value  = "2F Kids' Dorms, 2F Dorms Main Hall"
$rounds.Range("D9").Value  = "2F Kids' Dorms, 2F Dorms Main Hall"
$rounds.Range("D10").Value = "B Laundry Room, B Supply Room"

$rounds.Range("E83").Value = "Defense"

# ---------------------------------------------------------------------
# 2. Players sheet - correct mis-recorded "Spawn" site values
# ---------------------------------------------------------------------
$players = $wb.Worksheets.Item("Players")

$siteFixes = @{
    32 = "2F Kids' Dorms, 2F Dorms Main Hall"
    33 = "2F Kids' Dorms, 2F Dorms Main Hall"
    34 = "2F Kids' Dorms, 2F Dorms Main Hall"
    35 = "2F Kids' Dorms, 2F Dorms Main Hall"
    36 = "2F Kids' Dorms, 2F Dorms Main Hall"
    37 = "2F Kids' Dorms, 2F Dorms Main Hall"
    38 = "2F Kids' Dorms, 2F Dorms Main Hall"
    39 = "2F Kids' Dorms, 2F Dorms Main Hall"
    40 = "2F Kids' Dorms, 2F Dorms Main Hall"
    41 = "2F Kids' Dorms, 2F Dorms Main Hall"
    42 = "B Laundry Room, B Supply Room"
    43 = "B Laundry Room, B Supply Room"
    44 = "B Laundry Room, B Supply Room"
    45 = "B Laundry Room, B Supply Room"
    46 = "B Laundry Room, B Supply Room"
}

foreach ($row in $siteFixes.Keys) {
    $players.Cells.Item($row, 5).Value = $siteFixes[$row]
}

# ---------------------------------------------------------------------
# 3. Players sheet - normalize operator names to lowercase
# ---------------------------------------------------------------------
$operatorMap = @{
    "Iana"     = "iana"
    "Hibana"   = "hibana"
    "Sledge"   = "sledge"
    "Grim"     = "grim"
    "Ying"     = "ying"
    "Ace"      = "ace"
    "Thatcher" = "thatcher"
    "Maverick" = "maverick"
    "Ash"      = "ash"
    "Capitao"  = "capitao"
    "Amaru"    = "amaru"
    "Jager"    = "jäger"
    "Azami"    = "azami"
    "Fenrir"   = "fenrir"
    "Tachanka" = "tachanka"
    "Solis"    = "solis"
    "Wamai"    = "wamai"
    "Mute"     = "mute"
    "Castle"   = "castle"
    "Kaid"     = "kaid"
    "Frost"    = "frost"
    "Warden"   = "warden"
    "Kapkan"   = "kapkan"
    "Pulse"    = "pulse"
    "Mozzie"   = "mozzie"
    "Mira"     = "mira"
    "Aruni"    = "aruni"
    "Buck"     = "buck"
    "Nomad"    = "nomad"
    "Jackal"   = "jackal"
    "Zofia"    = "zofia"
    "Lion"     = "lion"
    "IQ"       = "iq"
    "Dokkaebi" = "dokkaebi"
    "Blitz"    = "blitz"
    "Bandit"   = "bandit"
    "Doc"      = "doc"
    "Finka"    = "finka"
    "Thermite" = "thermite"
    "Osa"      = "osa"
    "Twitch"   = "twitch"
    "Vigil"    = "vigil"
    "Goyo"     = "goyo"
    "Smoke"    = "smoke"
    "Maestro"  = "maestro"
    "Oryx"     = "oryx"
    "Brava"    = "brava"
    "Nokk"     = "nokk"
    "Gridlock" = "gridlock"
    "Echo"     = "echo"
    "Valkyrie" = "valkyrie"
    "Alibi"    = "alibi"
    "Fuze"     = "fuze"
}

$usedRows = $players.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $players.Cells.Item($r, 4)
    $current = $cell.Value2
    if ($operatorMap.ContainsKey($current)) {
        $cell.Value = $operatorMap[$current]
    }
}
